$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "B"="1.398111682538286"; "C"="1.102401346575959"; "D"="0.07823519779050514"; "E"="0.4092139868845805"; "G"="0.002463081803594892"; "I"="3.808639805109038" }
  3 = @{ "B"="1.282652077757803"; "C"="0.9851498860622314"; "D"="0.07093616729936514"; "E"="0.3561823458071416"; "G"="0.00247260389829089"; "I"="3.471471902225232" }
  4 = @{ "B"="1.21366884969143"; "C"="0.9139159272868937"; "D"="0.0665110432251339"; "E"="0.323841364384208"; "G"="0.002478722988872705"; "I"="3.265070426314281" }
  5 = @{ "B"="1.186024478400043"; "C"="0.8850687720020005"; "D"="0.06472143427211563"; "E"="0.3107122537241196"; "G"="0.002481285492114703"; "I"="3.181091008861102" }
  6 = @{ "B"="1.181461967248197"; "C"="0.8802893963844554"; "D"="0.06442508074619013"; "E"="0.3085350460549705"; "G"="0.002481715168467242"; "I"="3.167153464462956" }
  7 = @{ "B"="1.21329415503601"; "C"="0.9135261630189575"; "D"="0.0664868533134495"; "E"="0.3236641047673459"; "G"="0.002478757268086783"; "I"="3.263937353749583" }
  8 = @{ "B"="1.357896755383081"; "C"="1.061809661134475"; "D"="0.07570642896239121"; "E"="0.3908795828444624"; "G"="0.002466308713397018"; "I"="3.692239965218988" }
  9 = @{ "B"="1.657244715879415"; "C"="1.359071137040758"; "D"="0.09426047135147542"; "E"="0.5246942599270596"; "G"="0.002444040469492155"; "I"="4.538264888683386" }
  10 = @{ "B"="1.887722238917831"; "C"="1.582087147121456"; "D"="0.1082198027059462"; "E"="0.6246004464923658"; "G"="0.002428960074469515"; "I"="5.165316898411078" }
  11 = @{ "B"="1.995073650205313"; "C"="1.684697649966665"; "D"="0.114650155714429"; "E"="0.6704801522862169"; "G"="0.002422371745864983"; "I"="5.452161238919246" }
  12 = @{ "B"="2.036102833786344"; "C"="1.723732827116578"; "D"="0.1170974104021667"; "E"="0.6879227028873061"; "G"="0.0024199155378696"; "I"="5.561043209839113" }
  13 = @{ "B"="2.027249393530724"; "C"="1.71531775371659"; "D"="0.1165697956076883"; "E"="0.6841629622065426"; "G"="0.002420442814184889"; "I"="5.537581453435678" }
  14 = @{ "B"="1.998441478762629"; "C"="1.687905443133104"; "D"="0.1148512438921614"; "E"="0.671913740687927"; "G"="0.002422168900244034"; "I"="5.461113628399062" }
  15 = @{ "B"="1.980845486180669"; "C"="1.671138263429611"; "D"="0.1138001935771058"; "E"="0.6644199178644499"; "G"="0.002423231197697362"; "I"="5.414309700890556" }
  16 = @{ "B"="1.880758452117107"; "C"="1.575405662922549"; "D"="0.1078012358687346"; "E"="0.6216113658671958"; "G"="0.002429396070133433"; "I"="5.146605662365403" }
  17 = @{ "B"="1.820012161870238"; "C"="1.516982228229267"; "D"="0.1041420891959319"; "E"="0.5954650804879122"; "G"="0.002433247328778065"; "I"="4.982808141806629" }
  18 = @{ "B"="1.78530731529662"; "C"="1.483486993591441"; "D"="0.1020449493299083"; "E"="0.5804666580977766"; "G"="0.00243548808299335"; "I"="4.88874439055374" }
  19 = @{ "B"="1.773596667477761"; "C"="1.472164335721914"; "D"="0.1013361621087512"; "E"="0.575395156282724"; "G"="0.002436251175699488"; "I"="4.856920557902299" }
  20 = @{ "B"="1.8264542836597"; "C"="1.523190199515852"; "D"="0.1045308300015506"; "E"="0.5982441851489568"; "G"="0.00243283470834454"; "I"="5.000229103646319" }
  21 = @{ "B"="2.006892675328629"; "C"="1.695952152013888"; "D"="0.1153556869391537"; "E"="0.6755097086066399"; "G"="0.002421660861532926"; "I"="5.483566768715093" }
  22 = @{ "B"="2.12702860238312"; "C"="1.80991021978366"; "D"="0.1225019194875046"; "E"="0.7264121409753841"; "G"="0.002414583174080491"; "I"="5.800987778486217" }
  23 = @{ "B"="2.062701824610144"; "C"="1.748988641696997"; "D"="0.1186810584200231"; "E"="0.6992051910457775"; "G"="0.00241834022199519"; "I"="5.631423600737662" }
  24 = @{ "B"="1.823541119538277"; "C"="1.520383286383208"; "D"="0.104355059926533"; "E"="0.5969876483305967"; "G"="0.002433021171202395"; "I"="4.992352760122088" }
  25 = @{ "B"="1.574479944711584"; "C"="1.277892707884234"; "D"="0.08918641965993857"; "E"="0.4882416556510094"; "G"="0.00244983788456982"; "I"="4.308559645820083" }
}

foreach ($row in $data.Keys) {
  $rowData = $data[$row]
  foreach ($col in $rowData.Keys) {
    $ws.Range("$col$row").Value = [double]$rowData[$col]
  }
}
